$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10, shifting rows 10-11 down to 11-12
$ws.Rows.Item(10).Insert()

# Set the new cell value for the inserted slot
$ws.Range("A10").Value = "16:00 - 17:00"

# Update the selection
$ws.Range("A10").Select()

# Configure page setup (paper size 9 = A4, portrait orientation)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

